$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.825.68'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '2.919.46'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'357.37"
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').Value = "'109.42"
$ws.Range('E6').Value = '  -2.50%  '
$ws.Range('D7').Value = "'0.563"
$ws.Range('E7').Value = '  +1.10%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.626"
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').Value = "'39.08"
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = "'0.0872"
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.137"
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').Value = "'19.53"
$ws.Range('E13').Value = '  -1.97%  '
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '3.377.13'
$ws.Range('D16').Value = '2.919.54'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').Value = '51.803.48'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = "'3.36"
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('D20').Value = "'7.55"
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').Value = "'13.94"
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = "'70.60"
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = "'268.77"
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').Value = "'2.82"
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E26').Value = '  +13.13%  '
$ws.Range('D27').Value = "'7.75"
$ws.Range('E27').Value = '  +20.07%  '
$ws.Range('D28').Value = "'26.90"
$ws.Range('E28').Value = '  +0.72%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +9.41%  '
$ws.Range('D31').Value = "'10.52"
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').Value = "'37.35"
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').Value = "'2.21"
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = "'6.07"
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = "'52.27"
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('D36').Value = "'0.0443"
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = "'3.20"
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('D39').Value = "'18.25"
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('E41').Value = '  -3.78%  '
$ws.Range('E42').Value = '  +2.50%  '
$ws.Range('D43').Value = "'22.78"
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').Value = "'118.94"
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('E46').Value = '  -5.76%  '
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('D48').Value = '2.129.16'
$ws.Range('E48').Value = '  -3.16%  '
$ws.Range('D49').Value = "'0.249"
$ws.Range('E49').Value = '  -4.94%  '
$ws.Range('D50').Value = "'0.0338"
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').Value = "'9.11"
